# Historique.xlsx update
# - Remove the "Caméra 3ème personne" backlog task (row 27 cleared).
# - Add two new realised tasks in rows 22/23 ("Ajoute le sceptre ..." and
#   "Ajoute une caméra 3è personne"), each marked done ("OK") with a
#   "-" realised-in value, and durations 1h / 40min.
# - Fill the previously-empty "Réalisé en" column (C) for rows 12-15 with "-".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 12-15: column C ("Réalisé en") was blank, now "-"
$ws.Cells.Item(12, 3).Value = "-"
$ws.Cells.Item(13, 3).Value = "-"
$ws.Cells.Item(14, 3).Value = "-"
$ws.Cells.Item(15, 3).Value = "-"

# Row 22: new task "Ajoute le sceptre (avec effet de lumière, halo et particules)"
$ws.Cells.Item(22, 1).Value = "Ajoute le sceptre (avec effet de lumière, halo et particules)"
$ws.Cells.Item(22, 2).Value = "1h"
$ws.Cells.Item(22, 3).Value = "-"
$ws.Cells.Item(22, 6).Value = "OK"

# Row 23: new task "Ajoute une caméra 3è personne"
$ws.Cells.Item(23, 1).Value = "Ajoute une caméra 3è personne"
$ws.Cells.Item(23, 2).Value = "40min"
$ws.Cells.Item(23, 3).Value = "-"
$ws.Cells.Item(23, 6).Value = "OK"

# Row 27: clear the old "Caméra 3ème personne" backlog entry
$ws.Cells.Item(27, 1).ClearContents()
$ws.Cells.Item(27, 6).ClearContents()

# Update the view's selected cell to match the author's final cursor position.
$ws.Range("A19").Select()
